$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 - this shifts the existing row 35 (and
# everything below it) down by one, matching the diff where old rows
# 35..133 become new rows 36..134.
$ws.Rows.Item(35).Insert()

# Populate the newly-inserted row 35 with the new record.
# (Most columns mirror the record that used to sit at row 35, now at row 36,
# except for the date, volume, min/max/avg price and $/Kg columns.)
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 45054
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100108
$ws.Cells.Item(35, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(35, 9).Value = 100108002
$ws.Cells.Item(35, 10).Value = "Mango"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 60
$ws.Cells.Item(35, 14).Value = 8000
$ws.Cells.Item(35, 15).Value = 8000
$ws.Cells.Item(35, 16).Value = 8000
$ws.Cells.Item(35, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(35, 18).Value = "Perú"
$ws.Cells.Item(35, 19).Value = 2000
$ws.Cells.Item(35, 20).Value = 4
